$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 29: new translation entry "SingleUseId47" -> "FPS: <value>"
$ws.Cells.Item(29, 2).Value = "SingleUseId47"
$ws.Cells.Item(29, 3).Value = "Medium"
$ws.Cells.Item(29, 4).Value = "Left"
$ws.Cells.Item(29, 5).Value = "LTR"
$ws.Cells.Item(29, 6).Value = "FPS: <value>"

# Row 30: new translation entry "SingleUseId48" -> "0"
$ws.Cells.Item(30, 2).Value = "SingleUseId48"
$ws.Cells.Item(30, 3).Value = "Medium"
$ws.Cells.Item(30, 4).Value = "Left"
$ws.Cells.Item(30, 5).Value = "LTR"
$ws.Cells.Item(30, 6).NumberFormat = "@"
$ws.Cells.Item(30, 6).Value = "0"
$ws.Cells.Item(30, 6).Style = "Normal"
